# Generate Report for Archive
#
# The localization-status report lists each source file's row sorted by
# file name. Re-sorting moved "c2618bbd-...md" up from row 5 to row 3,
# pushing "fce21647-...md" and "923d649d-...md" down one row each (row 5
# keeps "923d649d-...md"). Apply that row re-sort, with its matching
# per-language handoff metadata and hyperlink display text, on all three
# sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: columns A (file), B (zh-cn status), C (de-de status)
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A3").Value = "c2618bbd-0c51-4380-8b70-b376139c604d.md"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"

$ws.Range("A4").Value = "fce21647-efde-4127-ac6e-9cee3f0f070f.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

$ws.Range("A5").Value = "923d649d-0171-4c9d-a1f2-32ac1819f07c.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "Ready for handoff"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "c2618bbd-0c51-4380-8b70-b376139c604d.md" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "fce21647-efde-4127-ac6e-9cee3f0f070f.md" }
    elseif ($addr -eq '$A$5') { $hl.TextToDisplay = "923d649d-0171-4c9d-a1f2-32ac1819f07c.md" }
}

# ---- zh-cn sheet: A=file, B=status, C=latest handoff file, D=latest handoff datetime
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A3").Value = "c2618bbd-0c51-4380-8b70-b376139c604d.md"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "c2618bbd-0c51-4380-8b70-b376139c604d.e595e65dfc7618a015615a210d25c4061ee8cbd7.zh-cn.xlf"
$ws.Range("D3").Value = "2016-02-17 05:51:07"

$ws.Range("A4").Value = "fce21647-efde-4127-ac6e-9cee3f0f070f.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "fce21647-efde-4127-ac6e-9cee3f0f070f.a1c10ac379cd3efd5bd676537b7e276c7e92ba86.zh-cn.xlf"
$ws.Range("D4").Value = "2016-02-17 05:49:41"

$ws.Range("A5").Value = "923d649d-0171-4c9d-a1f2-32ac1819f07c.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "923d649d-0171-4c9d-a1f2-32ac1819f07c.a807b1663a689ae6dd945a048415116076b59d8a.zh-cn.xlf"
$ws.Range("D5").Value = "2016-02-17 05:51:49"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "c2618bbd-0c51-4380-8b70-b376139c604d.md" }
    elseif ($addr -eq '$C$3') { $hl.TextToDisplay = "c2618bbd-0c51-4380-8b70-b376139c604d.e595e65dfc7618a015615a210d25c4061ee8cbd7.zh-cn.xlf" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "fce21647-efde-4127-ac6e-9cee3f0f070f.md" }
    elseif ($addr -eq '$C$4') { $hl.TextToDisplay = "fce21647-efde-4127-ac6e-9cee3f0f070f.a1c10ac379cd3efd5bd676537b7e276c7e92ba86.zh-cn.xlf" }
    elseif ($addr -eq '$A$5') { $hl.TextToDisplay = "923d649d-0171-4c9d-a1f2-32ac1819f07c.md" }
    elseif ($addr -eq '$C$5') { $hl.TextToDisplay = "923d649d-0171-4c9d-a1f2-32ac1819f07c.a807b1663a689ae6dd945a048415116076b59d8a.zh-cn.xlf" }
}

# ---- de-de sheet: A=file, B=status, C=latest handoff file, D=latest handoff datetime
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A3").Value = "c2618bbd-0c51-4380-8b70-b376139c604d.md"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "c2618bbd-0c51-4380-8b70-b376139c604d.e595e65dfc7618a015615a210d25c4061ee8cbd7.de-de.xlf"
$ws.Range("D3").Value = "2016-02-17 05:51:18"

$ws.Range("A4").Value = "fce21647-efde-4127-ac6e-9cee3f0f070f.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "fce21647-efde-4127-ac6e-9cee3f0f070f.a1c10ac379cd3efd5bd676537b7e276c7e92ba86.de-de.xlf"
$ws.Range("D4").Value = "2016-02-17 05:50:01"

$ws.Range("A5").Value = "923d649d-0171-4c9d-a1f2-32ac1819f07c.md"
$ws.Range("B5").Value = "Ready for handoff"
$ws.Range("C5").Value = "923d649d-0171-4c9d-a1f2-32ac1819f07c.a807b1663a689ae6dd945a048415116076b59d8a.de-de.xlf"
$ws.Range("D5").Value = "2016-02-17 05:51:59"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$3') { $hl.TextToDisplay = "c2618bbd-0c51-4380-8b70-b376139c604d.md" }
    elseif ($addr -eq '$C$3') { $hl.TextToDisplay = "c2618bbd-0c51-4380-8b70-b376139c604d.e595e65dfc7618a015615a210d25c4061ee8cbd7.de-de.xlf" }
    elseif ($addr -eq '$A$4') { $hl.TextToDisplay = "fce21647-efde-4127-ac6e-9cee3f0f070f.md" }
    elseif ($addr -eq '$C$4') { $hl.TextToDisplay = "fce21647-efde-4127-ac6e-9cee3f0f070f.a1c10ac379cd3efd5bd676537b7e276c7e92ba86.de-de.xlf" }
    elseif ($addr -eq '$A$5') { $hl.TextToDisplay = "923d649d-0171-4c9d-a1f2-32ac1819f07c.md" }
    elseif ($addr -eq '$C$5') { $hl.TextToDisplay = "923d649d-0171-4c9d-a1f2-32ac1819f07c.a807b1663a689ae6dd945a048415116076b59d8a.de-de.xlf" }
}
